$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Donation to own country
$ws.Range("B2").Value = 30.5598705501618
$ws.Range("C2").Value = 36.7049608355091
$ws.Range("D2").Value = 30.3550724637681
$ws.Range("E2").Value = 32.4707692307692
$ws.Range("F2").Value = 33.287775246773

# Row 3: Donation to Africa
$ws.Range("B3").Value = 32.2250803858521
$ws.Range("C3").Value = 40.3689839572192
$ws.Range("D3").Value = 30.0074906367041
$ws.Range("E3").Value = 33.1504702194357
$ws.Range("F3").Value = 30.7006033182504
